# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-obsolete detail rows for the old worker (NAREN, last period)
# and the extra periods for ANDREA that are no longer part of this statement.
# This leaves row16/row17 (NAREN's old rows, to be repointed to ANDREA below)
# and what was row22 (ANDREA's last row, which already carries the
# "closing row" border styling) shifted up to become the new row18.
$ws.Rows("18:21").Delete()

# Update the header summary figures.
$ws.Range("E11").Value = 170820
$ws.Range("C13").Value = 1
$ws.Range("F13").Value = 3

# Row 16: repoint from NAREN ANDRES CHICO RODRIGUEZ to ANDREA DEL CARMEN VITOLA AGUILERA,
# keep period 2506, and update the salary figure.
$ws.Range("C16").Value = "1143398929"
$ws.Range("D16").Value = "ANDREA DEL CARMEN VITOLA AGUILERA"
$ws.Range("E16").Value = "2506"
$ws.Range("G16").Value = 1423500

# Row 17: same worker, period 2507.
$ws.Range("C17").Value = "1143398929"
$ws.Range("D17").Value = "ANDREA DEL CARMEN VITOLA AGUILERA"
$ws.Range("E17").Value = "2507"
$ws.Range("G17").Value = 1423500

# Row 18 (formerly row22): same worker, new period 2508, full-month value.
$ws.Range("C18").Value = "1143398929"
$ws.Range("D18").Value = "ANDREA DEL CARMEN VITOLA AGUILERA"
$ws.Range("E18").Value = "2508"
$ws.Range("F18").Value = 56940
$ws.Range("G18").Value = 1423500
